# Update "predictions" column (column B) values for a set of rows, based on
# the new feature-prediction pass described in the commit message
# ("possible features to use").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new predicted value (0 or 1)
$updates = @{
    41  = 1
    46  = 0
    81  = 1
    115 = 1
    132 = 1
    134 = 0
    224 = 1
    258 = 1
    309 = 0
    315 = 1
    371 = 1
    407 = 0
    452 = 1
    462 = 1
    467 = 0
    562 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
